# Add data for 2021-11-22: the "through" date on the running November 2021
# column moves from Nov 13 to Nov 14, and the new day's carjacking counts
# are reflected in the B column (current month-to-date) plus a handful of
# historical cells that were corrected at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet/tab title and the matching header label in B1.
$ws.Name = "Through 2021-11-14"
$ws.Range("B1").Value = "November 2021 (through November 14)"

# Updated running-month (column B) counts by neighborhood row.
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 3
$ws.Range("B8").Value = 5
$ws.Range("B13").Value = 3
$ws.Range("B31").Value = 2
$ws.Range("B37").Value = 4
$ws.Range("B64").Value = 3

# Other corrected / newly populated historical cells.
$ws.Range("BP2").Value = 2
$ws.Range("M3").Value = 7
$ws.Range("AI3").Value = 3
$ws.Range("AT3").Value = 3
$ws.Range("BE5").Value = 4
$ws.Range("BP9").Value = 4
$ws.Range("M24").Value = 2
$ws.Range("M32").Value = 4
$ws.Range("M40").Value = 1
$ws.Range("AT40").Value = 1
$ws.Range("BE48").Value = 2
$ws.Range("BE62").Value = 1
$ws.Range("M80").Value = 1
$ws.Range("AT80").Value = 1
$ws.Range("X84").Value = 2
$ws.Range("M96").Value = 1
